$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    4   = -12.12970000000001
    7   = -13.1844
    16  = -13.9156
    28  = -12.9121
    29  = -11.2974
    32  = -13.1753
    40  = -13.0276
    52  = -11.2837
    57  = -13.80219999999999
    66  = -11.3936
    100 = -12.83679999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
